$d = $word.ActiveDocument

$d.Content.Find.Execute("Der Calliope ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Der Calliope ", 2)
